{"js": "const body = context.document.body;\n\nconst replacements = [\n  {\n    find: \"D:\\\\Career\\\\5_Project_to_add_to_CV\\\\1_House-price-analysis-with-BN\",\n    replace: \"D:\\\\Career\\\\5_Project_to_add_to_CV\\\\1_rent_analysis\"\n  },\n  {\n    find: \"git remote add origin https://github.com/username/repo-name.git\",\n    replace: \"git remote add origin https://github.com/taitran0102/rent-analysis.git\"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the clone path and the git remote URL throughout the document.\n$d = $word.ActiveDocument\n\n$oldPath = \"D:\\Career\\5_Project_to_add_to_CV\\1_House-price-analysis-with-BN\"\n$newPath = \"D:\\Career\\5_Project_to_add_to_CV\\1_rent_analysis\"\n\n$oldRemote = \"git remote add origin https://github.com/username/repo-name.git\"\n$newRemote = \"git remote add origin https://github.com/taitran0102/rent-analysis.git\"\n\n# Replace every occurrence of the old project path (appears 3 times in the doc).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldPath\n$find.Replacement.Text = $newPath\n$find.Execute(\n  $oldPath,\n  $false,\n  $false,\n  $false,\n  $false,\n  $false,\n  $true,\n  1,\n  $false,\n  $newPath,\n  2\n)\n\n# Replace the placeholder git remote URL with the real repo URL.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = $oldRemote\n$find2.Replacement.Text = $newRemote\n$find2.Execute(\n  $oldRemote,\n  $false,\n  $false,\n  $false,\n  $false,\n  $false,\n  $true,\n  1,\n  $false,\n  $newRemote,\n  2\n)\n"}
